$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 85, shifting existing rows 85-89 down to 86-90
$ws.Rows.Item(85).Insert()

# Fill in the new row 85 with a copy of the (now shifted) old-row-85 data (now in row 86),
# but with an updated date value (44931).
$ws.Cells.Item(85, 1).Value = 1
$ws.Cells.Item(85, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(85, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(85, 4).Value = 44931
$ws.Cells.Item(85, 5).Value = 15
$ws.Cells.Item(85, 6).Value = 100112012
$ws.Cells.Item(85, 7).Value = "Espinaca"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 250
$ws.Cells.Item(85, 11).Value = 2800
$ws.Cells.Item(85, 12).Value = 3000
$ws.Cells.Item(85, 13).Value = 2900
$ws.Cells.Item(85, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(85, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(85, 16).Value = 967
$ws.Cells.Item(85, 17).Value = 3
$ws.Cells.Item(85, 18).Value = "Hortaliza"

# Match the date cell style (style index 2, numFmt 165) used by other date cells in column D
$ws.Cells.Item(85, 4).NumberFormat = $ws.Cells.Item(86, 4).NumberFormat
